$d = $word.ActiveDocument

# 1. Replace the first sentence text (introduce typos / double space)
$d.Content.Find.Execute(
    "40. Hussein used a water collector to collect some water while on a trip to a desert.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "40.  fussein used a water cotlector to coliect some water while on a trip to a desert.",
    2)

# 2. Replace the (c) sentence text
$d.Content.Find.Execute(
    "(c) Explain why there are folds in the metal top.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    [char]0x0028 + "c) " + [char]0x201C + "Explain why there are folds in th" + [char]0x00E9 + " metal top.",
    2)

# 3. Remove the two paragraphs that contain the inline drawings (images)
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.InlineShapes.Count -gt 0) {
        $p.Range.Delete()
    }
}
